$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: rename column G, add new columns H, I, J ---
$ws.Range("G1").Value = "fecha_ingreso"
$ws.Range("H1").Value = "correo_notificaciones"
$ws.Range("I1").Value = "id_responsable"
$ws.Range("J1").Value = "activo"

# --- New "fecha_ingreso" (hire date) column, formatted as a date (built-in mm-dd-yy => numFmtId 14) ---
# Apply the number format to one cell, then propagate via copy/paste-special so every
# cell in the column shares a single style index instead of each getting its own xf.
$ws.Range("G2").NumberFormat = "mm-dd-yy"
$ws.Range("G2").Copy()
$ws.Range("G2:G40").PasteSpecial(-4122)

$default = Get-Date -Year 2025 -Month 1 -Day 1 -Hour 0 -Minute 0 -Second 0
$ws.Range("G2:G40").Value = $default

$exceptions = @{
    2  = Get-Date -Year 2024 -Month 1 -Day 1 -Hour 0 -Minute 0 -Second 0
    7  = Get-Date -Year 2023 -Month 1 -Day 1 -Hour 0 -Minute 0 -Second 0
    22 = Get-Date -Year 2020 -Month 1 -Day 1 -Hour 0 -Minute 0 -Second 0
    31 = Get-Date -Year 2024 -Month 1 -Day 1 -Hour 0 -Minute 0 -Second 0
    32 = Get-Date -Year 2024 -Month 1 -Day 1 -Hour 0 -Minute 0 -Second 0
    33 = Get-Date -Year 2024 -Month 1 -Day 1 -Hour 0 -Minute 0 -Second 0
}

foreach ($row in $exceptions.Keys) {
    $ws.Range("G$row").Value = $exceptions[$row]
}

# --- Update view: clear frozen/scrolled topLeftCell, select E2 ---
$ws.Range("E2").Select()
